$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M", now "B")
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.8571428571428571
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.9230769230769231
$ws.Range("E2").Value = 36

# Row 3 (was "B", now "M")
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.7142857142857143
$ws.Range("D3").Value = 0.8333333333333334
$ws.Range("E3").Value = 21

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8947368421052632
$ws.Range("C4").Value = 0.8947368421052632
$ws.Range("D4").Value = 0.8947368421052632
$ws.Range("E4").Value = 0.8947368421052632

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9285714285714286
$ws.Range("C5").Value = 0.8571428571428572
$ws.Range("D5").Value = 0.8782051282051282
$ws.Range("E5").Value = 57

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9097744360902255
$ws.Range("C6").Value = 0.8947368421052632
$ws.Range("D6").Value = 0.8900134952766532
$ws.Range("E6").Value = 57
